$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.310.38'

# Row 3
$ws.Range('D3').Value = '1.610.49'
$ws.Range('E3').Value = '  +0.53%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').Value = "'213.10"
$ws.Range('E5').Value = '  +0.35%  '

# Row 6
$ws.Range('E6').Value = '  -0.06%  '

# Row 7
$ws.Range('E7').Value = '  +0.33%  '

# Row 8
$ws.Range('E8').Value = '  +0.98%  '

# Row 9
$ws.Range('D9').Value = "'0.0615"
$ws.Range('E9').Value = '  +0.26%  '

# Row 10
$ws.Range('E10').Value = '  +2.63%  '

# Row 11
$ws.Range('E11').Value = '  -0.37%  '

# Row 12
$ws.Range('D12').Value = '1.835.18'

# Row 13
$ws.Range('D13').Value = '1.593.30'
$ws.Range('E13').Value = '  -0.55%  '

# Row 14
$ws.Range('E14').Value = '  +0.40%  '

# Row 15
$ws.Range('E15').Value = '  +1.05%  '

# Row 16
$ws.Range('D16').Value = '26.283.34'
$ws.Range('E16').Value = '  +0.81%  '

# Row 17
$ws.Range('D17').Value = "'62.28"
$ws.Range('E17').Value = '  +3.16%  '

# Row 18
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  +0.94%  '

# Row 19
$ws.Range('E19').Value = '  -0.08%  '

# Row 20
$ws.Range('D20').Value = "'201.94"
$ws.Range('E20').Value = '  +0.13%  '

# Row 21
$ws.Range('E21').Value = '  +1.20%  '

# Row 22
$ws.Range('D22').Value = "'9.34"
$ws.Range('E22').Value = '  +0.80%  '

# Row 23
$ws.Range('E23').Value = '  +0.90%  '

# Row 24
$ws.Range('E24').Value = '  +1.99%  '

# Row 25
$ws.Range('D25').Value = "'143.70"

# Row 26
$ws.Range('E26').Value = '  -0.06%  '

# Row 27
$ws.Range('E27').Value = '  -0.32%  '

# Row 28
$ws.Range('E28').Value = '  +0.65%  '

# Row 29
$ws.Range('D29').Value = "'6.57"
$ws.Range('E29').Value = '  +2.32%  '

# Row 30
$ws.Range('D30').Value = "'0.0497"
$ws.Range('E30').Value = '  +5.08%  '

# Row 31
$ws.Range('E31').Value = '  +0.17%  '

# Row 32
$ws.Range('E32').Value = '  +2.80%  '

# Row 33
$ws.Range('D33').Value = "'2.94"
$ws.Range('E33').Value = '  -0.35%  '

# Row 34
$ws.Range('E34').Value = '  +1.00%  '

# Row 35
$ws.Range('E35').Value = '  +1.15%  '

# Row 36
$ws.Range('D36').Value = '1.162.59'
$ws.Range('E36').Value = '  +3.37%  '

# Row 37
$ws.Range('E37').Value = '  +2.20%  '

# Row 38
$ws.Range('E38').Value = '  -0.08%  '

# Row 39
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = "'2.33"
$ws.Range('E39').Value = '  +1.04%  '

# Row 40
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = "'0.789"
$ws.Range('E40').Value = '  +0.07%  '

# Row 41
$ws.Range('D41').Value = "'0.496"
$ws.Range('E41').Value = '  +1.08%  '

# Row 42
$ws.Range('D42').Value = "'5.36"
$ws.Range('E42').Value = '  +4.10%  '

# Row 43
$ws.Range('E43').Value = '  +0.09%  '

# Row 44
$ws.Range('D44').Value = '1.746.14'
$ws.Range('E44').Value = '  +0.48%  '

# Row 45
$ws.Range('D45').Value = "'92.47"
$ws.Range('E45').Value = '  -0.60%  '

# Row 46
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0105'
$ws.Range('E46').Value = '  +14.01%  '

# Row 47
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = "'1.54"
$ws.Range('E47').Value = '  +1.83%  '

# Row 48
$ws.Range('D48').Value = "'53.90"
$ws.Range('E48').Value = '  +0.95%  '

# Row 49
$ws.Range('E49').Value = '  +0.61%  '

# Row 50
$ws.Range('E50').Value = '  -0.14%  '

# Row 51
$ws.Range('E51').Value = '  -0.12%  '

Write-Host "Applied cryptos update."